$d = $word.ActiveDocument

function Get-ParagraphIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r") -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "Create get_cars with IsAuthenticated"
#       -> "Create " / "cars_user" / " with IsAuthenticated"   (3 runs)
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Create get_cars with IsAuthenticated")
if ($found1) {
    $p1Start = $rng1.Start

    $wordRange = $d.Range($p1Start + 7, $p1Start + 7 + 8)
    $wordRange.Text = "cars_user"

    # Force the run to split at the word boundaries (same formatting either
    # side, so a plain text edit would otherwise get re-coalesced into a
    # single run).
    $splitA = $d.Range($p1Start + 7, $p1Start + 7)
    $d.Bookmarks.Add("zzz_split_a", $splitA)
    $splitB = $d.Range($p1Start + 7 + 9, $p1Start + 7 + 9)
    $d.Bookmarks.Add("zzz_split_b", $splitB)
    $d.Bookmarks.Item("zzz_split_a").Delete()
    $d.Bookmarks.Item("zzz_split_b").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Create create_car with IsAuthenticated"
#       -> "Update cars_user to accept POST request to create new car"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Create create_car with IsAuthenticated", $true, $false, $false, $false, $false, $true, 1, $false, "Update cars_user to accept POST request to create new car", 2)

# ---------------------------------------------------------------------------
# 3) Remove the now-redundant "Add path to Car's urls" bullet that used to
#    follow the create_car step.
# ---------------------------------------------------------------------------
if ($found2) {
    $updatedIdx = Get-ParagraphIndexByText("Update cars_user to accept POST request to create new car")
    if ($updatedIdx -gt 0 -and $updatedIdx -lt $d.Paragraphs.Count) {
        $nextPara = $d.Paragraphs.Item($updatedIdx + 1)
        if ($nextPara.Range.Text.TrimEnd("`r") -eq "Add path to Car’s urls") {
            $nextPara.Range.Delete()
        }
    }
}
